$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above the existing row 1069 (the block of
# "Femacal de La Calera" / Coliflor records is sorted with the most
# recent weekly price entries at the top of this sub-range), pushing
# the former rows 1069-1190 down to 1071-1192.
$insertRange = $ws.Range("A1069:R1070")
$insertRange.EntireRow.Insert()

# New row 1069: Primera quality entry for the new weekly date (45212 = 2023-10-13)
$ws.Range("A1069").Value = 3
$ws.Range("B1069").Value = "Femacal de La Calera"
$ws.Range("C1069").Value = "Coquimbo"
$ws.Range("D1069").Value = 45212
$ws.Range("D1069").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E1069").Value = 5
$ws.Range("F1069").Value = 100112008
$ws.Range("G1069").Value = "Coliflor"
$ws.Range("H1069").Value = "Sin especificar"
$ws.Range("I1069").Value = "Primera"
$ws.Range("J1069").Value = 2000
$ws.Range("K1069").Value = 750
$ws.Range("L1069").Value = 800
$ws.Range("M1069").Value = 772
$ws.Range("N1069").Value = "`$/unidad"
$ws.Range("O1069").Value = "Provincia de Quillota"
$ws.Range("P1069").Value = 772
$ws.Range("Q1069").Value = 1
$ws.Range("R1069").Value = "Hortaliza"

# New row 1070: Segunda quality entry for the same new weekly date
$ws.Range("A1070").Value = 3
$ws.Range("B1070").Value = "Femacal de La Calera"
$ws.Range("C1070").Value = "Coquimbo"
$ws.Range("D1070").Value = 45212
$ws.Range("D1070").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E1070").Value = 5
$ws.Range("F1070").Value = 100112008
$ws.Range("G1070").Value = "Coliflor"
$ws.Range("H1070").Value = "Sin especificar"
$ws.Range("I1070").Value = "Segunda"
$ws.Range("J1070").Value = 1500
$ws.Range("K1070").Value = 600
$ws.Range("L1070").Value = 600
$ws.Range("M1070").Value = 600
$ws.Range("N1070").Value = "`$/unidad"
$ws.Range("O1070").Value = "Provincia de Quillota"
$ws.Range("P1070").Value = 600
$ws.Range("Q1070").Value = 1
$ws.Range("R1070").Value = "Hortaliza"
